# Word COM-interop script applying the Arabic translation edits.
$d = $word.ActiveDocument

# wdReplaceAll = 2 ; wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $find"
    }
}

# 1 & 3: "English" -> "الإنجليزية" (both occurrences use the same translation)
Replace-Text "English" "الإنجليزية"

# 2: language list translated to Arabic
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

# 4: brief description partially translated
Replace-Text "An email sent to partners in the target country who have RSVPed no. It will be sent via customer.io" "An email sent to partners in the target country who have RSVPed no. سيتم إرسالها عبر customer.io"

# 6: trailing space removed
Replace-Text "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up. " "Even though we’re disappointed we can’t meet you, we understand that scheduling conflicts and other commitments sometimes come up."

# 7: trailing space removed
Replace-Text "We hope to see you at our future events. " "We hope to see you at our future events."

# 8: sentence translated to Arabic
Replace-Text "If you have any questions, please contact us via " "إذا كانت لديك أي أسئلة، فاتصل بنا:  "

# 9: hyperlink text translated
Replace-Text "live chat" "الدردشة الحية"

# 10: sentence translated to Arabic
Replace-Text "If you have any questions, please contact your country manager, " "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمدير بلدك  "

# 5: remove the extra empty run (plain space with rtl) at the end of the
#    "Subject line" paragraph - there were two trailing single-space runs,
#    keep only the highlighted one and delete the plain one.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -like "Subject line:*") {
        $runs = $p.Range.Words
        break
    }
}

$subjectPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Subject line")) {
        $subjectPara = $p
        break
    }
}
if ($subjectPara -ne $null) {
    $pEnd = $subjectPara.Range.End
    # paragraph mark is the last character; the last run ends right before it
    $lastRunRange = $d.Range($pEnd - 2, $pEnd - 1)
    if ($lastRunRange.Text -eq " ") {
        $lastRunRange.Delete()
    }
}

# 11: translate the comment text (Arabic)
$c = $d.Comments.Item(1)
$c.Text = "اختر أيًا منهما"
